$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Text edits to existing cells (shared-string content updates)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value  = "LLM Judge: Gemini 1.5"
$ws.Range("A9").Value  = "LLM Judge: Gemini 2.0 Experimental"
$ws.Range("B16").Value = "LLM to extract and classify claims; then # truthful claims / total # claims"
$ws.Range("B17").Value = "LLM to extract and classify claims; then # truthful claims / total # claims"
$ws.Range("C16").Value = "LLM request for each context node then weighted cumulative precision"
$ws.Range("C17").Value = "1 LLM request for all context nodes then weighted cumulative precision"
$ws.Range("B18").Value = "LLM; result calculated from PPI & validation set accuracy"
$ws.Range("C18").Value = "LLM; result calculated from PPI & validation set accuracy- “Context Relevance”"

# ---------------------------------------------------------------------------
# 2. Remove the trailing D15 cell (row 15 header no longer has a 4th column)
# ---------------------------------------------------------------------------
$ws.Range("D15").ClearContents()

# ---------------------------------------------------------------------------
# 3. Add new B28 value
# ---------------------------------------------------------------------------
$ws.Range("B28").Value = 0.745486278078553
$ws.Range("B28").NumberFormat = "0.000"

# ---------------------------------------------------------------------------
# 4. New block of rows (34-46): extra analysis tables
# ---------------------------------------------------------------------------

# Row 34 - section header (taller row)
$ws.Range("A34").Value = "Average Difference between LLM Results"
$ws.Range("B34").Value = "Faithfulness"
$ws.Range("C34").Value = "Contextual Precision"
$ws.Rows.Item(34).RowHeight = 28

# Row 35
$ws.Range("A35").Value = "RAGAS"
$ws.Range("B35").Value = 0.077
$ws.Range("C35").Value = 0.069

# Row 36
$ws.Range("A36").Value = "DeepEval"
$ws.Range("B36").Value = 0.054
$ws.Range("C36").Value = 0.145

# Row 38 - sub header
$ws.Range("B38").Value = "Faithfulness"
$ws.Range("C38").Value = "Contextual Precision"

# Row 39 - section header
$ws.Range("A39").Value = "Average Difference between Examples"

# Row 40
$ws.Range("A40").Value = "RAGAS and DeepEval, Gemini 1.5"
$ws.Range("B40").Value = 0.161216071534399
$ws.Range("C40").Value = 0.172519308304842
$ws.Range("B40").NumberFormat = "0.000"
$ws.Range("C40").NumberFormat = "0.000"

# Row 41
$ws.Range("A41").Value = "RAGAS and DeepEval, Gemini 2.0 Exp"
$ws.Range("B41").Value = 0.122093926909796
$ws.Range("C41").Value = 0.183165000802064
$ws.Range("B41").NumberFormat = "0.000"
$ws.Range("C41").NumberFormat = "0.000"

# Row 43 - sub header
$ws.Range("B43").Value = "Faithfulness"
$ws.Range("C43").Value = "Contextual Precision"

# Row 44 - section header
$ws.Range("A44").Value = "% of Examples within 0.1 scores of each other"

# Row 45
$ws.Range("A45").Value = "RAGAS and DeepEval, Gemini 1.5"
$ws.Range("B45").Formula = "=640/903"
$ws.Range("C45").Formula = "=525/903"
$ws.Range("B45").NumberFormat = "0.00%"
$ws.Range("C45").NumberFormat = "0.00%"

# Row 46
$ws.Range("A46").Value = "RAGAS and DeepEval, Gemini 2.0 Exp"
$ws.Range("B46").Formula = "=709/903"
$ws.Range("C46").Formula = "=506/903"
$ws.Range("B46").NumberFormat = "0.00%"
$ws.Range("C46").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# 5. Column widths (characters) and view/selection state
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 37.666666666666664
$ws.Columns.Item(2).ColumnWidth = 58.333333333333336
$ws.Columns.Item(3).ColumnWidth = 17.833333333333332

$ws.Range("B13").Select() | Out-Null
